$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 202
$ws.Range("I2").Value = 169.11111
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 169.11111
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -56.11111
$ws.Range("N2").Value = -576
$ws.Range("H46").Value = 7490
$ws.Range("J46").Value = 6653.3335
$ws.Range("L46").Value = 19960.0005
$ws.Range("N46").Value = -20198.0005
$ws.Range("H60").Value = 7490
$ws.Range("J60").Value = 6653.3335
$ws.Range("L60").Value = 19960.0005
$ws.Range("N60").Value = -20928.0005
$ws.Range("H76").Value = 13892302
$ws.Range("I76").Value = 111111110
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 111111110
$ws.Range("L76").Value = 3900
$ws.Range("M76").Value = -111110795
$ws.Range("N76").Value = -4530
$ws.Range("H79").Value = 13892302
$ws.Range("I79").Value = 111111110
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 111111110
$ws.Range("L79").Value = 3900
$ws.Range("M79").Value = -111110018
$ws.Range("N79").Value = -6084
$ws.Range("H129").Value = 1121.85
$ws.Range("J129").Value = 1160.1052
$ws.Range("L129").Value = 3480.3156
$ws.Range("N129").Value = -13480.3156
$ws.Range("H132").Value = 26130.195
$ws.Range("I132").Value = 26130.195
$ws.Range("K132").Value = 78390.58499999999
$ws.Range("M132").Value = -75860.58499999999
$ws.Range("H137").Value = 27028474
$ws.Range("I137").Value = 37038064
$ws.Range("J137").Value = 2588.3
$ws.Range("K137").Value = 111114192
$ws.Range("L137").Value = 7764.900000000001
$ws.Range("M137").Value = -111111642
$ws.Range("N137").Value = -12864.9
$ws.Range("H138").Value = 6394419
$ws.Range("I138").Value = 1691200
$ws.Range("J138").Value = 9806558
$ws.Range("K138").Value = 5073600
$ws.Range("L138").Value = 29419674
$ws.Range("M138").Value = -5068460
$ws.Range("N138").Value = -29429954

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 60189.293
$ws.Range("I2").Value = 63888.625
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 63888.625
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -63775.625
$ws.Range("N2").Value = -1226
$ws.Range("H63").Value = 10057
$ws.Range("I63").Value = 10931.75
$ws.Range("J63").Value = 9557.143
$ws.Range("K63").Value = 10931.75
$ws.Range("L63").Value = 9557.143
$ws.Range("M63").Value = -10245.75
$ws.Range("N63").Value = -10929.143
$ws.Range("H66").Value = 10057
$ws.Range("I66").Value = 10931.75
$ws.Range("J66").Value = 9557.143
$ws.Range("K66").Value = 54658.75
$ws.Range("L66").Value = 47785.715
$ws.Range("M66").Value = -51226.75
$ws.Range("N66").Value = -54649.715
$ws.Range("H116").Value = 60189.293
$ws.Range("I116").Value = 63888.625
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 63888.625
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = -61594.625
$ws.Range("N116").Value = -5588
$ws.Range("H122").Value = 10543.667
$ws.Range("I122").Value = 12052.4
$ws.Range("K122").Value = 36157.2
$ws.Range("M122").Value = -33707.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 60189.293
$ws.Range("I3").Value = 63888.625
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 63888.625
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -63774.625
$ws.Range("N3").Value = -1228
$ws.Range("H7").Value = 3334833.8
$ws.Range("I7").Value = 1750.75
$ws.Range("J7").Value = 10001000
$ws.Range("K7").Value = 1750.75
$ws.Range("L7").Value = 10001000
$ws.Range("M7").Value = -1637.75
$ws.Range("N7").Value = -10001226
$ws.Range("H105").Value = 246963.95
$ws.Range("I105").Value = 2835.7144
$ws.Range("K105").Value = 2835.7144
$ws.Range("M105").Value = -1088.7144
$ws.Range("H126").Value = 57250
$ws.Range("J126").Value = 58500
$ws.Range("L126").Value = 58500
$ws.Range("N126").Value = -68380
$ws.Range("H134").Value = 3684.0356
$ws.Range("I134").Value = 2681.476
$ws.Range("K134").Value = 8044.428
$ws.Range("M134").Value = -5509.428

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4839.137
$ws.Range("I31").Value = 1311.0278
$ws.Range("J31").Value = 13306.6
$ws.Range("K31").Value = 1311.0278
$ws.Range("L31").Value = 13306.6
$ws.Range("M31").Value = -1016.0278
$ws.Range("N31").Value = -13896.6
$ws.Range("H34").Value = 4839.137
$ws.Range("I34").Value = 1311.0278
$ws.Range("J34").Value = 13306.6
$ws.Range("K34").Value = 1311.0278
$ws.Range("L34").Value = 13306.6
$ws.Range("M34").Value = -1109.0278
$ws.Range("N34").Value = -13710.6
$ws.Range("H58").Value = 1911.1
$ws.Range("I58").Value = 1148.75
$ws.Range("J58").Value = 4960.5
$ws.Range("K58").Value = 1148.75
$ws.Range("L58").Value = 4960.5
$ws.Range("M58").Value = -945.75
$ws.Range("N58").Value = -5366.5
$ws.Range("H105").Value = 870.625
$ws.Range("I105").Value = 870.625
$ws.Range("K105").Value = 870.625
$ws.Range("M105").Value = 876.375
$ws.Range("H132").Value = 1545.4359
$ws.Range("I132").Value = 927.8182
$ws.Range("J132").Value = 4942.3335
$ws.Range("K132").Value = 2783.4546
$ws.Range("L132").Value = 14827.0005
$ws.Range("M132").Value = -253.4546
$ws.Range("N132").Value = -19887.0005
$ws.Range("H136").Value = 1911.1
$ws.Range("I136").Value = 1148.75
$ws.Range("J136").Value = 4960.5
$ws.Range("K136").Value = 3446.25
$ws.Range("L136").Value = 14881.5
$ws.Range("M136").Value = -896.25
$ws.Range("N136").Value = -19981.5
$ws.Range("H141").Value = 532127.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 532127.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 532127.2
$ws.Range("N141").Value = -542487.2
$ws.Range("M141").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4244.3335
$ws.Range("I56").Value = 4244.3335
$ws.Range("K56").Value = 4244.3335
$ws.Range("M56").Value = -3714.3335
$ws.Range("H131").Value = 4017378.8
$ws.Range("I131").Value = 666.6667
$ws.Range("J131").Value = 4168005.5
$ws.Range("K131").Value = 2000.0001
$ws.Range("L131").Value = 12504016.5
$ws.Range("M131").Value = 3039.9999
$ws.Range("N131").Value = -12514096.5
$ws.Range("H132").Value = 1389.4348
$ws.Range("I132").Value = 1227.4286
$ws.Range("J132").Value = 1460.3125
$ws.Range("K132").Value = 11046.8574
$ws.Range("L132").Value = 13142.8125
$ws.Range("M132").Value = -8516.857399999999
$ws.Range("N132").Value = -18202.8125
$ws.Range("H136").Value = 2851.4773
$ws.Range("I136").Value = 1806
$ws.Range("J136").Value = 2985.513
$ws.Range("K136").Value = 5418
$ws.Range("L136").Value = 8956.539000000001
$ws.Range("M136").Value = -318
$ws.Range("N136").Value = -19156.539
$ws.Range("H140").Value = 8364.093999999999
$ws.Range("I140").Value = 13888.75
$ws.Range("K140").Value = 41666.25
$ws.Range("M140").Value = -36486.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 45002.5
$ws.Range("I19").Value = 45002.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 45002.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -44714.5
$ws.Range("H102").Value = 2447.8235
$ws.Range("I102").Value = 2368.4285
$ws.Range("J102").Value = 2576.077
$ws.Range("K102").Value = 2368.4285
$ws.Range("L102").Value = 2576.077
$ws.Range("M102").Value = -746.4285
$ws.Range("N102").Value = -5820.077
$ws.Range("N19").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2878.3215
$ws.Range("I7").Value = 2443.2856
$ws.Range("J7").Value = 3023.3333
$ws.Range("K7").Value = 2443.2856
$ws.Range("L7").Value = 3023.3333
$ws.Range("M7").Value = -2331.2856
$ws.Range("N7").Value = -3247.3333
$ws.Range("H17").Value = 2786.6667
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3680
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 3680
$ws.Range("M17").Value = -830
$ws.Range("N17").Value = -4020
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H22").Value = 1320.1428
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 1256.8334
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 1256.8334
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -1846.8334
$ws.Range("H27").Value = 1320.1428
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 1256.8334
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 1256.8334
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -1470.8334
$ws.Range("H126").Value = 2878.3215
$ws.Range("I126").Value = 2443.2856
$ws.Range("J126").Value = 3023.3333
$ws.Range("K126").Value = 7329.8568
$ws.Range("L126").Value = 9069.999899999999
$ws.Range("M126").Value = -4859.8568
$ws.Range("N126").Value = -14009.9999
$ws.Range("M18").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1211.3334
$ws.Range("I122").Value = 1043.1428
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 3129.4284
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -679.4284000000002
$ws.Range("N122").Value = -10300
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
